## Update KNIVSTA overview: reorder data rows 2-117 according to the
## new sort order captured in $rowMap (new row -> old row), and bump
## the "Förändrad" (changed) date in column C from 46076 to 46077 for
## every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 26
$firstRow = 2
$lastRow = 117
$newDateSerial = 46077

# Map of new-row -> old-row (both are 1-based Excel COM row indices,
# matching the worksheet's r="N" attribute in the underlying XML).
$rowMap = @{}
$rowMap[2] = 2
$rowMap[3] = 3
$rowMap[4] = 5
$rowMap[5] = 4
$rowMap[6] = 6
$rowMap[7] = 7
$rowMap[8] = 8
$rowMap[9] = 17
$rowMap[10] = 13
$rowMap[11] = 14
$rowMap[12] = 19
$rowMap[13] = 18
$rowMap[14] = 10
$rowMap[15] = 9
$rowMap[16] = 11
$rowMap[17] = 12
$rowMap[18] = 15
$rowMap[19] = 16
$rowMap[20] = 20
$rowMap[21] = 21
$rowMap[22] = 22
$rowMap[23] = 23
$rowMap[24] = 24
$rowMap[25] = 25
$rowMap[26] = 26
$rowMap[27] = 28
$rowMap[28] = 29
$rowMap[29] = 30
$rowMap[30] = 31
$rowMap[31] = 33
$rowMap[32] = 32
$rowMap[33] = 36
$rowMap[34] = 34
$rowMap[35] = 35
$rowMap[36] = 27
$rowMap[37] = 37
$rowMap[38] = 38
$rowMap[39] = 39
$rowMap[40] = 40
$rowMap[41] = 41
$rowMap[42] = 74
$rowMap[43] = 69
$rowMap[44] = 108
$rowMap[45] = 109
$rowMap[46] = 47
$rowMap[47] = 95
$rowMap[48] = 86
$rowMap[49] = 76
$rowMap[50] = 48
$rowMap[51] = 88
$rowMap[52] = 111
$rowMap[53] = 106
$rowMap[54] = 71
$rowMap[55] = 46
$rowMap[56] = 107
$rowMap[57] = 105
$rowMap[58] = 57
$rowMap[59] = 102
$rowMap[60] = 87
$rowMap[61] = 49
$rowMap[62] = 110
$rowMap[63] = 45
$rowMap[64] = 50
$rowMap[65] = 112
$rowMap[66] = 91
$rowMap[67] = 53
$rowMap[68] = 101
$rowMap[69] = 104
$rowMap[70] = 42
$rowMap[71] = 61
$rowMap[72] = 63
$rowMap[73] = 64
$rowMap[74] = 44
$rowMap[75] = 52
$rowMap[76] = 68
$rowMap[77] = 70
$rowMap[78] = 97
$rowMap[79] = 98
$rowMap[80] = 96
$rowMap[81] = 54
$rowMap[82] = 73
$rowMap[83] = 72
$rowMap[84] = 75
$rowMap[85] = 43
$rowMap[86] = 77
$rowMap[87] = 113
$rowMap[88] = 114
$rowMap[89] = 115
$rowMap[90] = 78
$rowMap[91] = 79
$rowMap[92] = 82
$rowMap[93] = 81
$rowMap[94] = 80
$rowMap[95] = 83
$rowMap[96] = 85
$rowMap[97] = 84
$rowMap[98] = 116
$rowMap[99] = 117
$rowMap[100] = 89
$rowMap[101] = 90
$rowMap[102] = 92
$rowMap[103] = 94
$rowMap[104] = 93
$rowMap[105] = 51
$rowMap[106] = 100
$rowMap[107] = 103
$rowMap[108] = 99
$rowMap[109] = 58
$rowMap[110] = 56
$rowMap[111] = 55
$rowMap[112] = 59
$rowMap[113] = 60
$rowMap[114] = 62
$rowMap[115] = 67
$rowMap[116] = 65
$rowMap[117] = 66

# --- 1. Snapshot every existing data cell (value or formula) ---------
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $entry = @{}
        if ($cell.HasFormula) {
            $entry["hasFormula"] = $true
            $entry["formula"] = $cell.Formula
        } else {
            $entry["hasFormula"] = $false
            $entry["value"] = $cell.Value2
        }
        $rowData += $entry
    }
    $snapshot[$r] = $rowData
}

# --- 2. Write the snapshot back out in the new row order -------------
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $rowData = $snapshot[$oldRow]
    for ($c = 1; $c -le $numCols; $c++) {
        $entry = $rowData[$c - 1]
        $cell = $ws.Cells.Item($newRow, $c)
        if ($entry["hasFormula"]) {
            $cell.Formula = $entry["formula"]
        } else {
            $val = $entry["value"]
            if ($null -eq $val) {
                $cell.Value = ""
            } else {
                $cell.Value = $val
            }
        }
    }
    # Column C ("Förändrad") always advances to the new stamp date.
    $ws.Cells.Item($newRow, 3).Value = $newDateSerial
}

# --- 3. Re-pin every data row to its original 15pt height -------------
# Writing into the wrap-text "Artnamn" column (R, style s=2) can make the
# engine auto-fit taller rows; restore the authored 15pt height so the
# sheet's row metrics stay exactly as before (matches the source file,
# where every row was an explicit ht="15" customHeight="1").
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}
